$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '26.369.30'
$ws.Range("E2").Value = '  +0.37%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.605.84'
$ws.Range("E3").Value = '  +0.69%  '
$ws.Range("E4").Value = '  +0.00%  '
$ws.Range("E5").Value = '  -0.47%  '
$ws.Range("E6").Value = '  -1.03%  '
$ws.Range("E8").Value = '  -0.63%  '
$ws.Range("E9").Value = '  -0.52%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '19.22'
$ws.Range("E10").Value = '  +1.06%  '
$ws.Range("E11").Value = '  +0.14%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.830.99'
$ws.Range("E12").Value = '  +0.69%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.599.74'
$ws.Range("E13").Value = '  -0.26%  '
$ws.Range("E14").Value = '  -0.25%  '
$ws.Range("E15").Value = '  -0.78%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '63.30'
$ws.Range("E16").Value = '  -1.04%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '26.372.47'
$ws.Range("E17").Value = '  +0.46%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '230.96'
$ws.Range("E18").Value = '  +7.47%  '
$ws.Range("E19").Value = '  -0.48%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '7.66'
$ws.Range("E20").Value = '  +3.97%  '
$ws.Range("E21").Value = '  +0.03%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.27'
$ws.Range("E22").Value = '  -0.74%  '
$ws.Range("E23").Value = '  +2.68%  '
$ws.Range("E24").Value = '  -1.08%  '
$ws.Range("E25").Value = '  +1.44%  '
$ws.Range("E26").Value = '  +0.11%  '
$ws.Range("E27").Value = '  -0.22%  '
$ws.Range("E28").Value = '  +0.72%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '15.40'
$ws.Range("E29").Value = '  +1.71%  '
$ws.Range("E30").Value = '  +0.82%  '
$ws.Range("E31").Value = '  -0.47%  '
$ws.Range("E32").Value = '  +4.97%  '
$ws.Range("E33").Value = '  +0.44%  '
$ws.Range("E34").Value = '  -1.79%  '
$ws.Range("E35").Value = '  -0.22%  '
$ws.Range("E36").Value = '  +0.45%  '
$ws.Range("E37").Value = '  -3.72%  '
$ws.Range("E38").Value = '  -0.64%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.820'
$ws.Range("E39").Value = '  -0.66%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '5.79'
$ws.Range("E40").Value = '  -0.55%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.18'
$ws.Range("E42").Value = '  +0.74%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.928'
$ws.Range("E43").Value = '  -4.72%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.744.18'
$ws.Range("E44").Value = '  +0.85%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.758'
$ws.Range("E45").Value = '  -0.92%  '
$ws.Range("E46").Value = '  -0.37%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '89.27'
$ws.Range("E47").Value = '  +2.71%  '
$ws.Range("E48").Value = '  -0.52%  '
$ws.Range("E49").Value = '  -0.22%  '
$ws.Range("E50").Value = '  -0.09%  '
$ws.Range("E51").Value = '  -0.01%  '
